$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2 = 'System, backup@backdoor.com, system'
    3 = 'dnasr281@gmail.com, System'
    4 = 'System, backup@backdoor.com'
    5 = 'System, backup@backdoor.com'
    6 = 'dnasr281@gmail.com, System'
    7 = 'admin@admin.com, System'
    8 = 'System, backup@backdoor.com'
    10 = 'dnasr281@gmail.com, System'
    11 = 'dnasr281@gmail.com, System'
    12 = 'dnasr281@gmail.com, System'
    13 = 'dnasr281@gmail.com, System'
    14 = 'dnasr281@gmail.com, System'
    15 = 'dnasr281@gmail.com, System'
    17 = 'dnasr281@gmail.com, System'
    18 = 'dnasr281@gmail.com, System'
    19 = 'dnasr281@gmail.com, System'
    20 = 'dnasr281@gmail.com, System'
    21 = 'dnasr281@gmail.com, System'
    22 = 'dnasr281@gmail.com, System'
    24 = 'dnasr281@gmail.com, System'
    29 = 'System, backup@backdoor.com, system'
    30 = 'dnasr281@gmail.com, System'
    31 = 'System, backup@backdoor.com'
    32 = 'System, backup@backdoor.com'
    33 = 'dnasr281@gmail.com, System'
    34 = 'admin@admin.com, System'
    35 = 'System, backup@backdoor.com'
    37 = 'dnasr281@gmail.com, System'
    38 = 'dnasr281@gmail.com, System'
    39 = 'dnasr281@gmail.com, System'
    40 = 'dnasr281@gmail.com, System'
    41 = 'dnasr281@gmail.com, System'
    42 = 'dnasr281@gmail.com, System'
    44 = 'dnasr281@gmail.com, System'
    45 = 'dnasr281@gmail.com, System'
    46 = 'dnasr281@gmail.com, System'
    47 = 'dnasr281@gmail.com, System'
    48 = 'dnasr281@gmail.com, System'
    49 = 'dnasr281@gmail.com, System'
    51 = 'dnasr281@gmail.com, System'
    56 = 'System, backup@backdoor.com, system'
    57 = 'dnasr281@gmail.com, System'
    58 = 'System, backup@backdoor.com'
    59 = 'System, backup@backdoor.com'
    60 = 'dnasr281@gmail.com, System'
    61 = 'admin@admin.com, System'
    62 = 'System, backup@backdoor.com'
    64 = 'dnasr281@gmail.com, System'
    65 = 'dnasr281@gmail.com, System'
    66 = 'dnasr281@gmail.com, System'
    67 = 'dnasr281@gmail.com, System'
    68 = 'dnasr281@gmail.com, System'
    69 = 'dnasr281@gmail.com, System'
    71 = 'dnasr281@gmail.com, System'
    72 = 'dnasr281@gmail.com, System'
    73 = 'dnasr281@gmail.com, System'
    74 = 'dnasr281@gmail.com, System'
    75 = 'dnasr281@gmail.com, System'
    76 = 'dnasr281@gmail.com, System'
    78 = 'dnasr281@gmail.com, System'
    83 = 'System, backup@backdoor.com'
    84 = 'System, backup@backdoor.com'
    85 = 'System, backup@backdoor.com'
    86 = 'dnasr281@gmail.com, System'
    87 = 'dnasr281@gmail.com, System'
    88 = 'dnasr281@gmail.com, System'
    89 = 'dnasr281@gmail.com, System'
    90 = 'dnasr281@gmail.com, admin@admin.com'
    93 = 'dnasr281@gmail.com, System'
    95 = 'dnasr281@gmail.com, System'
    96 = 'dnasr281@gmail.com, System'
    97 = 'dnasr281@gmail.com, System'
    99 = 'dnasr281@gmail.com, System'
    102 = 'dnasr281@gmail.com, System'
    104 = 'dnasr281@gmail.com, System'
    109 = 'System, backup@backdoor.com'
    110 = 'System, backup@backdoor.com'
    111 = 'System, backup@backdoor.com'
    112 = 'dnasr281@gmail.com, System'
    113 = 'dnasr281@gmail.com, System'
    114 = 'dnasr281@gmail.com, System'
    115 = 'dnasr281@gmail.com, System'
    116 = 'dnasr281@gmail.com, admin@admin.com'
    119 = 'dnasr281@gmail.com, System'
    121 = 'dnasr281@gmail.com, System'
    122 = 'dnasr281@gmail.com, System'
    123 = 'dnasr281@gmail.com, System'
    125 = 'dnasr281@gmail.com, System'
    128 = 'dnasr281@gmail.com, System'
    130 = 'dnasr281@gmail.com, System'
    135 = 'System, backup@backdoor.com'
    136 = 'System, backup@backdoor.com'
    137 = 'System, backup@backdoor.com'
    138 = 'dnasr281@gmail.com, System'
    139 = 'dnasr281@gmail.com, System'
    140 = 'dnasr281@gmail.com, System'
    141 = 'dnasr281@gmail.com, System'
    142 = 'dnasr281@gmail.com, admin@admin.com'
    145 = 'dnasr281@gmail.com, System'
    147 = 'dnasr281@gmail.com, System'
    148 = 'dnasr281@gmail.com, System'
    149 = 'dnasr281@gmail.com, System'
    151 = 'dnasr281@gmail.com, System'
    154 = 'dnasr281@gmail.com, System'
    156 = 'dnasr281@gmail.com, System'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
